$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 2.288150666666667
$ws.Range("H2").Value = 6.864452
$ws.Range("I2").Value = 0.3964219041944151
$ws.Range("J2").Value = 0.3964219041944151
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 5.436778333333334
$ws.Range("N2").Value = 16.310335
$ws.Range("O2").Value = 0.121853993972124
$ws.Range("P2").Value = 0.121853993972124
$ws.Range("Q2").Value = 12.44016796793556
$ws.Range("R2").Value = 111.96151171142
$ws.Range("S2").Value = 0.04830559232412419
$ws.Range("T2").Value = 0.04830559232412419

# Row 3
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 2.288150666666667
$ws.Range("H3").Value = 6.864452
$ws.Range("I3").Value = 0.3964219041944151
$ws.Range("J3").Value = 0.3964219041944151
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 20.81735866666667
$ws.Range("N3").Value = 62.452076
$ws.Range("O3").Value = 0.4665774732677552
$ws.Range("P3").Value = 0.4665774732677551
$ws.Range("Q3").Value = 47.63325311137245
$ws.Range("R3").Value = 428.699278002352
$ws.Range("S3").Value = 0.1849615304070223
$ws.Range("T3").Value = 0.1849615304070223

# Row 4
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 2.288150666666667
$ws.Range("H4").Value = 6.864452
$ws.Range("I4").Value = 0.3964219041944151
$ws.Range("J4").Value = 0.3964219041944151
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.131792
$ws.Range("N4").Value = 15.395376
$ws.Range("O4").Value = 0.1150183643869107
$ws.Range("P4").Value = 0.1150183643869107
$ws.Range("Q4").Value = 11.74231328599467
$ws.Range("R4").Value = 105.680819573952
$ws.Range("S4").Value = 0.04559579902758623
$ws.Range("T4").Value = 0.04559579902758623

# Row 5
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 2.288150666666667
$ws.Range("H5").Value = 6.864452
$ws.Range("I5").Value = 0.3964219041944151
$ws.Range("J5").Value = 0.3964219041944151
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 13.23122433333333
$ws.Range("N5").Value = 39.693673
$ws.Range("O5").Value = 0.2965501683732101
$ws.Range("P5").Value = 0.2965501683732101
$ws.Range("Q5").Value = 30.27503477913288
$ws.Range("R5").Value = 272.475313012196
$ws.Range("S5").Value = 0.1175589824356823
$ws.Range("T5").Value = 0.1175589824356823

# Row 6
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 2.588894
$ws.Range("H6").Value = 7.766681999999999
$ws.Range("I6").Value = 0.4485256605643812
$ws.Range("J6").Value = 0.4485256605643813
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 5.436778333333334
$ws.Range("N6").Value = 16.310335
$ws.Range("O6").Value = 0.121853993972124
$ws.Range("P6").Value = 0.121853993972124
$ws.Range("Q6").Value = 14.07524280649667
$ws.Range("R6").Value = 126.67718525847
$ws.Range("S6").Value = 0.05465464313875505
$ws.Range("T6").Value = 0.05465464313875507

# Row 7
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.588894
$ws.Range("H7").Value = 7.766681999999999
$ws.Range("I7").Value = 0.4485256605643812
$ws.Range("J7").Value = 0.4485256605643813
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 20.81735866666667
$ws.Range("N7").Value = 62.452076
$ws.Range("O7").Value = 0.4665774732677552
$ws.Range("P7").Value = 0.4665774732677551
$ws.Range("Q7").Value = 53.89393494798133
$ws.Range("R7").Value = 485.045414531832
$ws.Range("S7").Value = 0.2092719694018798
$ws.Range("T7").Value = 0.2092719694018798

# Row 8
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.588894
$ws.Range("H8").Value = 7.766681999999999
$ws.Range("I8").Value = 0.4485256605643812
$ws.Range("J8").Value = 0.4485256605643813
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.131792
$ws.Range("N8").Value = 15.395376
$ws.Range("O8").Value = 0.1150183643869107
$ws.Range("P8").Value = 0.1150183643869107
$ws.Range("Q8").Value = 13.285665518048
$ws.Range("R8").Value = 119.570989662432
$ws.Range("S8").Value = 0.05158868786367381
$ws.Range("T8").Value = 0.05158868786367382

# Row 9
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.588894
$ws.Range("H9").Value = 7.766681999999999
$ws.Range("I9").Value = 0.4485256605643812
$ws.Range("J9").Value = 0.4485256605643813
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 13.23122433333333
$ws.Range("N9").Value = 39.693673
$ws.Range("O9").Value = 0.2965501683732101
$ws.Range("P9").Value = 0.2965501683732101
$ws.Range("Q9").Value = 34.25423728922066
$ws.Range("R9").Value = 308.288135602986
$ws.Range("S9").Value = 0.1330103601600725
$ws.Range("T9").Value = 0.1330103601600725

# Row 10
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 0.3123523333333333
$ws.Range("H10").Value = 0.9370569999999999
$ws.Range("I10").Value = 0.05411501461132016
$ws.Range("J10").Value = 0.05411501461132018
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 5.436778333333334
$ws.Range("N10").Value = 16.310335
$ws.Range("O10").Value = 0.121853993972124
$ws.Range("P10").Value = 0.121853993972124
$ws.Range("Q10").Value = 1.698190398232778
$ws.Range("R10").Value = 15.283713584095
$ws.Range("S10").Value = 0.006594130664249211
$ws.Range("T10").Value = 0.006594130664249213

# Row 11
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 0.3123523333333333
$ws.Range("H11").Value = 0.9370569999999999
$ws.Range("I11").Value = 0.05411501461132016
$ws.Range("J11").Value = 0.05411501461132018
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 20.81735866666667
$ws.Range("N11").Value = 62.452076
$ws.Range("O11").Value = 0.4665774732677552
$ws.Range("P11").Value = 0.4665774732677551
$ws.Range("Q11").Value = 6.502350553370221
$ws.Range("R11").Value = 58.52115498033199
$ws.Range("S11").Value = 0.02524884678319742
$ws.Range("T11").Value = 0.02524884678319742

# Row 12
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 0.3123523333333333
$ws.Range("H12").Value = 0.9370569999999999
$ws.Range("I12").Value = 0.05411501461132016
$ws.Range("J12").Value = 0.05411501461132018
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.131792
$ws.Range("N12").Value = 15.395376
$ws.Range("O12").Value = 0.1150183643869107
$ws.Range("P12").Value = 0.1150183643869107
$ws.Range("Q12").Value = 1.602927205381333
$ws.Range("R12").Value = 14.426344848432
$ws.Range("S12").Value = 0.006224220469367818
$ws.Range("T12").Value = 0.00622422046936782

# Row 13
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 0.3123523333333333
$ws.Range("H13").Value = 0.9370569999999999
$ws.Range("I13").Value = 0.05411501461132016
$ws.Range("J13").Value = 0.05411501461132018
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 13.23122433333333
$ws.Range("N13").Value = 39.693673
$ws.Range("O13").Value = 0.2965501683732101
$ws.Range("P13").Value = 0.2965501683732101
$ws.Range("Q13").Value = 4.132803793373443
$ws.Range("R13").Value = 37.19523414036099
$ws.Range("S13").Value = 0.01604781669450572
$ws.Range("T13").Value = 0.01604781669450572

# Row 14
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 0.5826116666666666
$ws.Range("H14").Value = 1.747835
$ws.Range("I14").Value = 0.1009374206298835
$ws.Range("J14").Value = 0.1009374206298836
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 5.436778333333334
$ws.Range("N14").Value = 16.310335
$ws.Range("O14").Value = 0.121853993972124
$ws.Range("P14").Value = 0.121853993972124
$ws.Range("Q14").Value = 3.167530486080556
$ws.Range("R14").Value = 28.507774374725
$ws.Range("S14").Value = 0.01229962784499558
$ws.Range("T14").Value = 0.01229962784499558

# Row 15
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 0.5826116666666666
$ws.Range("H15").Value = 1.747835
$ws.Range("I15").Value = 0.1009374206298835
$ws.Range("J15").Value = 0.1009374206298836
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 20.81735866666667
$ws.Range("N15").Value = 62.452076
$ws.Range("O15").Value = 0.4665774732677552
$ws.Range("P15").Value = 0.4665774732677551
$ws.Range("Q15").Value = 12.12843602838444
$ws.Range("R15").Value = 109.15592425546
$ws.Range("S15").Value = 0.04709512667565565
$ws.Range("T15").Value = 0.04709512667565565

# Row 16
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 0.5826116666666666
$ws.Range("H16").Value = 1.747835
$ws.Range("I16").Value = 0.1009374206298835
$ws.Range("J16").Value = 0.1009374206298836
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.131792
$ws.Range("N16").Value = 15.395376
$ws.Range("O16").Value = 0.1150183643869107
$ws.Range("P16").Value = 0.1150183643869107
$ws.Range("Q16").Value = 2.989841890106666
$ws.Range("R16").Value = 26.90857701096
$ws.Range("S16").Value = 0.01160965702628282
$ws.Range("T16").Value = 0.01160965702628282

# Row 17
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 0.5826116666666666
$ws.Range("H17").Value = 1.747835
$ws.Range("I17").Value = 0.1009374206298835
$ws.Range("J17").Value = 0.1009374206298836
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 13.23122433333333
$ws.Range("N17").Value = 39.693673
$ws.Range("O17").Value = 0.2965501683732101
$ws.Range("P17").Value = 0.2965501683732101
$ws.Range("Q17").Value = 7.708665660883888
$ws.Range("R17").Value = 69.37799094795498
$ws.Range("S17").Value = 0.02993300908294949
$ws.Range("T17").Value = 0.0299330090829495

Write-Output "Applied Inhba-Acvr1 natmi update (Dr Hou advice)"